$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 13: "~" / "Tilde" punctuator entry, appended after the existing
# table (rows 1-12 are the header + 11 existing punctuator rows).
$ws.Range("A13").Value = "~"
$ws.Range("B13").Value = "Tilde"

# A13 styling: Arial 10, bordered (thin, black), centered horizontally
# (matches the other data rows' border+font but without vertical centering).
$ws.Range("A13").Font.Name = "Arial"
$ws.Range("A13").Font.Size = 10
$ws.Range("A13").Borders.LineStyle = 1
$ws.Range("A13").Borders.Color = 0
$ws.Range("A13").HorizontalAlignment = -4108

# B13 styling: Arial 10, bordered (thin, black), default (left) alignment.
$ws.Range("B13").Font.Name = "Arial"
$ws.Range("B13").Font.Size = 10
$ws.Range("B13").Borders.LineStyle = 1
$ws.Range("B13").Borders.Color = 0

# Match the author's final cursor position/selection.
[void]$ws.Range("F11").Select()
